$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 header updates ---

# F1: "Walking" -> "Walking/swimming/flying"
$ws.Range("F1").Value = "Walking/swimming/flying"

# O1:T1 "Salvacion:*" columns get reordered
$ws.Range("O1").Value = "Salvacion:STR"
$ws.Range("P1").Value = "Salvacion:DEX"
$ws.Range("Q1").Value = "Salvacion:CON"
$ws.Range("R1").Value = "Salvacion:INT"
$ws.Range("S1").Value = "Salvacion:WIS"
$ws.Range("T1").Value = "Salvacion:CHA"

# U1:AM1 new skill columns (replacing the old bare skill names, each now annotated
# with its governing ability and a couple of passive-perception variants collapsed
# into one, plus the "Raligion" typo fixed to "Religion")
$ws.Range("U1").Value  = "Passive WIS (percepción)"
$ws.Range("V1").Value  = "Acrobatict (DEX)"
$ws.Range("W1").Value  = "Animal Handling (WIS)"
$ws.Range("X1").Value  = "Arcana (INT)"
$ws.Range("Y1").Value  = "Athletics (STR)"
$ws.Range("Z1").Value  = "Deception (CHA)"
$ws.Range("AA1").Value = "History (INT)"
$ws.Range("AB1").Value = "Insight (WIS)"
$ws.Range("AC1").Value = "Intimidation (CHA)"
$ws.Range("AD1").Value = "Investigation (INT)"
$ws.Range("AE1").Value = "Medicine (WIS)"
$ws.Range("AF1").Value = "Nature (INT)"
$ws.Range("AG1").Value = "Perception (WIS)"
$ws.Range("AH1").Value = "Performance (CHA)"
$ws.Range("AI1").Value = "Persusasion (CHA)"
$ws.Range("AJ1").Value = "Religion (INT)"
$ws.Range("AK1").Value = "Sleight of hand (DEX)"
$ws.Range("AL1").Value = "Stealth (DEX)"
$ws.Range("AM1").Value = "Survival (WIS)"

# New header block (F1:AM1) is centered
$ws.Range("F1:AM1").HorizontalAlignment = -4108

# --- Selection moves from F12 to F4 ---
[void]$ws.Range("F4").Select()
